$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1901234567901235
$ws.Range("C2").Value = 0.5753086419753086
$ws.Range("J2").Value = 0.009876543209876543
$ws.Range("P2").Value = 0.1234567901234568
$ws.Range("S2").Value = 0.1012345679012346

$ws.Range("B3").Value = 0.01666666666666667
$ws.Range("C3").Value = 0.02916666666666667
$ws.Range("J3").Value = 0.04166666666666666
$ws.Range("P3").Value = 0.7291666666666666
$ws.Range("S3").Value = 0.1833333333333333

$ws.Range("J4").Value = 0.02777777777777778
$ws.Range("P4").Value = 0.6388888888888888
$ws.Range("S4").Value = 0.3333333333333333

$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.5

$ws.Range("B6").Value = 0.06324110671936758
$ws.Range("D6").Value = 0.003952569169960474
$ws.Range("F6").Value = 0.02371541501976284
$ws.Range("J6").Value = 0.3359683794466403
$ws.Range("O6").Value = 0.01976284584980237
$ws.Range("Q6").Value = 0.1185770750988142
$ws.Range("R6").Value = 0.1067193675889328
$ws.Range("S6").Value = 0.3280632411067194

$ws.Range("B7").Value = 0.1116504854368932
$ws.Range("D7").Value = 0.02912621359223301
$ws.Range("F7").Value = 0.04854368932038835
$ws.Range("J7").Value = 0.145631067961165
$ws.Range("O7").Value = 0.01941747572815534
$ws.Range("Q7").Value = 0.2087378640776699
$ws.Range("R7").Value = 0.04854368932038835
$ws.Range("S7").Value = 0.3883495145631068

$ws.Range("B8").Value = 0.1447963800904978
$ws.Range("D8").Value = 0.01131221719457014
$ws.Range("E8").Value = 0.004524886877828055
$ws.Range("F8").Value = 0.07692307692307693
$ws.Range("J8").Value = 0.1108597285067873
$ws.Range("O8").Value = 0.02262443438914027
$ws.Range("Q8").Value = 0.1561085972850679
$ws.Range("R8").Value = 0.07466063348416289
$ws.Range("S8").Value = 0.3981900452488688

$ws.Range("B9").Value = 0.0847457627118644
$ws.Range("D9").Value = 0.005649717514124294
$ws.Range("F9").Value = 0.05649717514124294
$ws.Range("J9").Value = 0.06779661016949153
$ws.Range("O9").Value = 0.02259887005649718
$ws.Range("Q9").Value = 0.2542372881355932
$ws.Range("R9").Value = 0.07909604519774012
$ws.Range("S9").Value = 0.4293785310734463

$ws.Range("B10").Value = 0.1389854065323141
$ws.Range("D10").Value = 0.01598332175121612
$ws.Range("F10").Value = 0.07018763029881862
$ws.Range("J10").Value = 0.13273106323836
$ws.Range("O10").Value = 0.01945795691452398
$ws.Range("Q10").Value = 0.1980542043085476
$ws.Range("R10").Value = 0.08200138985406533
$ws.Range("S10").Value = 0.3425990271021543

$ws.Range("G11").Value = 0.1477572559366755
$ws.Range("J11").Value = 0.1451187335092348
$ws.Range("K11").Value = 0.237467018469657
$ws.Range("L11").Value = 0.4591029023746702
$ws.Range("S11").Value = 0.01055408970976253

$ws.Range("G12").Value = 0.6722222222222223
$ws.Range("K12").Value = 0.005555555555555556
$ws.Range("L12").Value = 0.03333333333333333
$ws.Range("S12").Value = 0.03888888888888889

$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.2592592592592592
$ws.Range("S13").Value = 0.07407407407407407

$ws.Range("F15").Value = 0.02904564315352697
$ws.Range("H15").Value = 0.1535269709543569
$ws.Range("I15").Value = 0.06639004149377593
$ws.Range("J15").Value = 0.3692946058091287
$ws.Range("K15").Value = 0.08298755186721991
$ws.Range("M15").Value = 0.01659751037344398
$ws.Range("O15").Value = 0.07468879668049792
$ws.Range("S15").Value = 0.2074688796680498

$ws.Range("F16").Value = 0.025
$ws.Range("I16").Value = 0.075
$ws.Range("J16").Value = 0.3875
$ws.Range("K16").Value = 0.125
$ws.Range("M16").Value = 0.02083333333333333
$ws.Range("O16").Value = 0.07083333333333333
$ws.Range("S16").Value = 0.1708333333333333

$ws.Range("F17").Value = 0.0170940170940171
$ws.Range("H17").Value = 0.1538461538461539
$ws.Range("I17").Value = 0.08547008547008547
$ws.Range("J17").Value = 0.4145299145299146
$ws.Range("K17").Value = 0.1431623931623932
$ws.Range("M17").Value = 0.0235042735042735
$ws.Range("N17").Value = 0.002136752136752137
$ws.Range("O17").Value = 0.0576923076923077
$ws.Range("S17").Value = 0.1025641025641026

$ws.Range("F18").Value = 0.04455445544554455
$ws.Range("H18").Value = 0.1881188118811881
$ws.Range("I18").Value = 0.06930693069306931
$ws.Range("J18").Value = 0.4207920792079208
$ws.Range("K18").Value = 0.1138613861386139
$ws.Range("M18").Value = 0.009900990099009901
$ws.Range("O18").Value = 0.06435643564356436
$ws.Range("S18").Value = 0.0891089108910891

$ws.Range("F19").Value = 0.02681992337164751
$ws.Range("H19").Value = 0.2030651340996169
$ws.Range("I19").Value = 0.0681992337164751
$ws.Range("J19").Value = 0.3846743295019157
$ws.Range("K19").Value = 0.1118773946360153
$ws.Range("M19").Value = 0.02681992337164751
$ws.Range("O19").Value = 0.06896551724137931
$ws.Range("S19").Value = 0.1095785440613027
